# Insert a new column before column A, shifting the existing data
# (old A:E) right to (new B:F), then populate the new column A with
# an "ID" header and per-row identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; existing columns shift right.
$ws.Columns.Item(1).Insert()

# Header for the new ID column. Copy the formatting already used by the
# other header cells (bold font, thin box border, centered) so the new
# A1 cell shares the same style entry as B1:F1, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ID"

# Row identifiers for the new column A (rows 2-25).
$ids = @(
    "Hb 2",
    "Hb 3",
    "S 24",
    "S 28",
    "Hb 107",
    "Hb 66",
    "Hb 69",
    "Hb 95",
    "Hb 99",
    "Hb 92",
    "Hb 40",
    "Hb 41",
    "S 11",
    "Hb 57",
    "S 21",
    "S 22",
    "S 3",
    "S 4",
    "S 5",
    "Hb 74",
    "Hb 79",
    "Hb 32",
    "S 15",
    "S 16"
)

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
